# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets.
# These two sheets mirror the same event listing, so the same cell updates
# are applied to both.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 12835
    "F6"  = 67
    "F12" = 13698
    "F13" = 14148
    "F22" = 1074
    "F26" = 5250
    "F28" = 280
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
